# Auto-generated edit script: updates crafting profit values across ALC, ARM, BSM, CRP, GSM, LTW, WVR sheets
# per scheduled market-data refresh (see commit message).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (index 1) ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("H43").Value = 1016.73334
$ws.Range("J43").Value = 1105
$ws.Range("L43").Value = 1105
$ws.Range("N43").Value = -1243
$ws.Range("H59").Value = 804.25
$ws.Range("I59").Value = 717
$ws.Range("J59").Value = 833.3333
$ws.Range("K59").Value = 2151
$ws.Range("L59").Value = 2499.9999
$ws.Range("M59").Value = -1594
$ws.Range("N59").Value = -3613.9999
$ws.Range("H64").Value = 3988.087
$ws.Range("I64").Value = 3545.375
$ws.Range("K64").Value = 3545.375
$ws.Range("M64").Value = -3297.375
$ws.Range("H67").Value = 3988.087
$ws.Range("I67").Value = 3545.375
$ws.Range("K67").Value = 3545.375
$ws.Range("M67").Value = -2687.375
$ws.Range("H74").Value = 4208.5713
$ws.Range("I74").Value = 3825
$ws.Range("J74").Value = 4720
$ws.Range("K74").Value = 3825
$ws.Range("L74").Value = 4720
$ws.Range("M74").Value = -2889
$ws.Range("N74").Value = -6592
$ws.Range("H77").Value = 4208.5713
$ws.Range("I77").Value = 3825
$ws.Range("J77").Value = 4720
$ws.Range("K77").Value = 19125
$ws.Range("L77").Value = 23600
$ws.Range("M77").Value = -14445
$ws.Range("N77").Value = -32960
$ws.Range("H129").Value = 861.7143
$ws.Range("J129").Value = 916.44183
$ws.Range("L129").Value = 2749.32549
$ws.Range("N129").Value = -12749.32549
$ws.Range("H137").Value = 1526.1562
$ws.Range("I137").Value = 1068.9642
$ws.Range("K137").Value = 3206.8926
$ws.Range("M137").Value = -656.8925999999997
$ws.Range("H138").Value = 1823.0172
$ws.Range("J138").Value = 3524.3076
$ws.Range("L138").Value = 10572.9228
$ws.Range("N138").Value = -20852.9228

# ---- Sheet: ARM (index 2) ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("H61").Value = 985.3333
$ws.Range("I61").Value = 882.4
$ws.Range("J61").Value = 1500
$ws.Range("K61").Value = 882.4
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = -670.4
$ws.Range("N61").Value = -1924
$ws.Range("H102").Value = 1450
$ws.Range("I102").Value = 1450
$ws.Range("K102").Value = 1450
$ws.Range("M102").Value = 172
$ws.Range("H136").Value = 985.3333
$ws.Range("I136").Value = 882.4
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 2647.2
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -97.19999999999982
$ws.Range("N136").Value = -9600

# ---- Sheet: BSM (index 3) ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("H15").Value = 50000
$ws.Range("J15").Value = 50000
$ws.Range("L15").Value = 50000
$ws.Range("N15").Value = -50454
$ws.Range("H86").Value = 2649.9167
$ws.Range("I86").Value = 3112.1667
$ws.Range("J86").Value = 2187.6667
$ws.Range("K86").Value = 3112.1667
$ws.Range("L86").Value = 2187.6667
$ws.Range("M86").Value = -1989.1667
$ws.Range("N86").Value = -4433.6667
$ws.Range("H89").Value = 2649.9167
$ws.Range("I89").Value = 3112.1667
$ws.Range("J89").Value = 2187.6667
$ws.Range("K89").Value = 15560.8335
$ws.Range("L89").Value = 10938.3335
$ws.Range("M89").Value = -9944.833500000001
$ws.Range("N89").Value = -22170.3335

# ---- Sheet: CRP (index 4) ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("H16").Value = 2152.5833
$ws.Range("I16").Value = 2232
$ws.Range("J16").Value = 1755.5
$ws.Range("K16").Value = 2232
$ws.Range("L16").Value = 1755.5
$ws.Range("M16").Value = -1945
$ws.Range("N16").Value = -2329.5
$ws.Range("H113").Value = 2152.5833
$ws.Range("I113").Value = 2232
$ws.Range("J113").Value = 1755.5
$ws.Range("K113").Value = 2232
$ws.Range("L113").Value = 1755.5
$ws.Range("M113").Value = -62
$ws.Range("N113").Value = -6095.5
$ws.Range("H132").Value = 2863.7036
$ws.Range("I132").Value = 1153.3334
$ws.Range("J132").Value = 5001.6665
$ws.Range("K132").Value = 3460.0002
$ws.Range("L132").Value = 15004.9995
$ws.Range("M132").Value = -930.0001999999999
$ws.Range("N132").Value = -20064.9995
$ws.Range("H134").Value = 1657.7222
$ws.Range("I134").Value = 1657.7222
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4973.1666
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2438.1666
$ws.Range("N134").ClearContents()

# ---- Sheet: GSM (index 6) ----
$ws = $wb.Worksheets.Item(6)
$ws.Range("H70").Value = 5611.1943
$ws.Range("I70").Value = 4840
$ws.Range("J70").Value = 7363.909
$ws.Range("K70").Value = 4840
$ws.Range("L70").Value = 7363.909
$ws.Range("M70").Value = -4570
$ws.Range("N70").Value = -7903.909
$ws.Range("H73").Value = 5611.1943
$ws.Range("I73").Value = 4840
$ws.Range("J73").Value = 7363.909
$ws.Range("K73").Value = 4840
$ws.Range("L73").Value = 7363.909
$ws.Range("M73").Value = -3904
$ws.Range("N73").Value = -9235.909
$ws.Range("H132").Value = 5187.1
$ws.Range("I132").Value = 5274.7144
$ws.Range("J132").Value = 4982.6665
$ws.Range("K132").Value = 15824.1432
$ws.Range("L132").Value = 14947.9995
$ws.Range("M132").Value = -13294.1432
$ws.Range("N132").Value = -20007.9995

# ---- Sheet: LTW (index 7) ----
$ws = $wb.Worksheets.Item(7)
$ws.Range("H16").Value = 1619.625
$ws.Range("I16").Value = 977.5
$ws.Range("J16").Value = 1833.6666
$ws.Range("K16").Value = 977.5
$ws.Range("L16").Value = 1833.6666
$ws.Range("M16").Value = -807.5
$ws.Range("N16").Value = -2173.6666
$ws.Range("H22").Value = 684.37036
$ws.Range("I22").Value = 374.16666
$ws.Range("J22").Value = 932.5333000000001
$ws.Range("K22").Value = 374.16666
$ws.Range("L22").Value = 932.5333000000001
$ws.Range("M22").Value = -79.16665999999998
$ws.Range("N22").Value = -1522.5333
$ws.Range("H27").Value = 684.37036
$ws.Range("I27").Value = 374.16666
$ws.Range("J27").Value = 932.5333000000001
$ws.Range("K27").Value = 374.16666
$ws.Range("L27").Value = 932.5333000000001
$ws.Range("M27").Value = -267.16666
$ws.Range("N27").Value = -1146.5333
$ws.Range("H46").Value = 146143
$ws.Range("I46").Value = 337000.34
$ws.Range("K46").Value = 337000.34
$ws.Range("M46").Value = -336812.34
$ws.Range("H55").Value = 496.91177
$ws.Range("I55").Value = 416.92593
$ws.Range("K55").Value = 416.92593
$ws.Range("M55").Value = -243.92593
$ws.Range("H139").Value = 37180.086
$ws.Range("J139").Value = 37961
$ws.Range("L139").Value = 37961
$ws.Range("N139").Value = -48241

# ---- Sheet: WVR (index 8) ----
$ws = $wb.Worksheets.Item(8)
$ws.Range("H62").Value = 3191.6667
$ws.Range("J62").Value = 3191.6667
$ws.Range("L62").Value = 3191.6667
$ws.Range("N62").Value = -4439.6667
$ws.Range("H65").Value = 3191.6667
$ws.Range("J65").Value = 3191.6667
$ws.Range("L65").Value = 15958.3335
$ws.Range("N65").Value = -22198.3335
$ws.Range("H132").Value = 1296.8679
$ws.Range("I132").Value = 989.27905
$ws.Range("J132").Value = 2619.5
$ws.Range("K132").Value = 2967.83715
$ws.Range("L132").Value = 7858.5
$ws.Range("M132").Value = -437.8371499999998
$ws.Range("N132").Value = -12918.5
$ws.Range("H136").Value = 1529.75
$ws.Range("I136").Value = 1142.8684
$ws.Range("J136").Value = 2999.9
$ws.Range("K136").Value = 3428.6052
$ws.Range("L136").Value = 8999.700000000001
$ws.Range("M136").Value = -878.6052
$ws.Range("N136").Value = -14099.7
